$d = $word.ActiveDocument
$sec = $d.Sections(1)
$hdr = $sec.Headers(1)
$result = $hdr.Range.Find.Execute("8:15/0", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)
Write-Output "Find/Replace result: $result"
